$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.966.78'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '''2.791.86'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''363.51'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').Value = '''109.92'
$ws.Range('E6').Value = '  -3.54%  '
$ws.Range('E7').Value = '  -2.19%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').Value = '''40.18'
$ws.Range('E10').Value = '  -3.64%  '
$ws.Range('D11').Value = '''0.0849'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').Value = '''19.50'
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = '''7.57'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').Value = '''3.224.32'
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('D16').Value = '''2.801.32'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').Value = '''51.915.22'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').Value = '''7.49'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').Value = '''13.15'
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('D22').Value = '''0.0₃0977'
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('D23').Value = '''70.41'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '''270.04'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '''2.77'
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').Value = '''26.56'
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '''0.162'
$ws.Range('E28').Value = '  +15.59%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''10.29'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.26'
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = '''51.98'
$ws.Range('E32').Value = '  -3.56%  '
$ws.Range('D33').Value = '''34.19'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '''5.75'
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('D35').Value = '''0.0846'
$ws.Range('D36').Value = '''5.25'
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('E39').Value = '  -2.27%  '
$ws.Range('D40').Value = '''2.00'
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''22.05'
$ws.Range('E44').Value = '  -7.55%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').Value = '''119.76'
$ws.Range('E45').Value = '  -6.70%  '
$ws.Range('D46').Value = '''2.085.60'
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('E47').Value = '  -4.17%  '
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('E50').Value = '  -6.08%  '
$ws.Range('D51').Value = '''8.87'
$ws.Range('E51').Value = '  -2.19%  '
